$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.217.74"
$ws.Range("E2").Value = "  +1.37%  "

# Row 3
$ws.Range("D3").Value = "1.819.32"
$ws.Range("E3").Value = "  -1.89%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.59"
$ws.Range("E5").Value = "  +2.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  +0.78%  "

# Row 7
$ws.Range("E7").Value = "  +0.32%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.80"
$ws.Range("E8").Value = "  -0.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.328"
$ws.Range("E9").Value = "  +7.97%  "

# Row 10
$ws.Range("E10").Value = "  -0.69%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.100"
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("D12").Value = "2.082.61"
$ws.Range("E12").Value = "  -1.56%  "

# Row 13
$ws.Range("D13").Value = "1.829.19"
$ws.Range("E13").Value = "  -1.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.14"
$ws.Range("E14").Value = "  -2.05%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.662"
$ws.Range("E15").Value = "  +0.77%  "

# Row 16
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.68"
$ws.Range("E16").Value = "  -0.69%  "

# Row 17
$ws.Range("D17").Value = "35.140.77"
$ws.Range("E17").Value = "  +1.21%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.83"
$ws.Range("E18").Value = "  +1.45%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0795"
$ws.Range("E19").Value = "  +0.90%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.27"
$ws.Range("E20").Value = "  -1.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.97"
$ws.Range("E21").Value = "  -1.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.67"
$ws.Range("E22").Value = "  -2.61%  "

# Row 23
$ws.Range("E23").Value = "  +0.38%  "

# Row 24
$ws.Range("E24").Value = "  +3.29%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.76"
$ws.Range("E25").Value = "  +0.27%  "

# Row 26
$ws.Range("E26").Value = "  -0.48%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.55"
$ws.Range("E27").Value = "  -1.84%  "

# Row 28
$ws.Range("E28").Value = "  -1.09%  "

# Row 29
$ws.Range("E29").Value = "  +24.08%  "

# Row 30
$ws.Range("E30").Value = "  +0.19%  "

# Row 31
$ws.Range("D31").Value = "3.334.58"
$ws.Range("E31").Value = "  +37.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.05"
$ws.Range("E32").Value = "  +2.78%  "

# Row 33
$ws.Range("E33").Value = "  +5.25%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.01"
$ws.Range("E34").Value = "  +0.17%  "

# Row 35
$ws.Range("E35").Value = "  -5.65%  "

# Row 36
$ws.Range("B36").Value = "Aave"
$ws.Range("C36").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "93.58"
$ws.Range("E36").Value = "  +2.58%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.15"
$ws.Range("E37").Value = "  +7.35%  "

# Row 38
$ws.Range("E38").Value = "  +2.67%  "

# Row 39
$ws.Range("E39").Value = "  +0.69%  "

# Row 40
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.29"
$ws.Range("E40").Value = "  +3.03%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.312.85"
$ws.Range("E41").Value = "  -2.40%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -2.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.78"
$ws.Range("E43").Value = "  -0.65%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.33"
$ws.Range("E44").Value = "  -4.09%  "

# Row 45
$ws.Range("E45").Value = "  +0.88%  "

# Row 46
$ws.Range("E46").Value = "  -2.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.35"
$ws.Range("E47").Value = "  +5.16%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  -1.32%  "

# Row 49
$ws.Range("D49").Value = "1.998.33"
$ws.Range("E49").Value = "  -0.59%  "

# Row 50
$ws.Range("E50").Value = "  +0.22%  "

# Row 51
$ws.Range("E51").Value = "  +5.48%  "
